$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: was Slotboard/298/98 -> now Bar-line/297/150.5
$ws.Range("A3").Value = "Bar-line"
$ws.Range("B3").Value = 297
$ws.Range("C3").Value = 150.5

# Row 4: was Bar-line / "297 px" / "150.5 px" -> now Bet Buttons/30/30
$ws.Range("A4").Value = "Bet Buttons"
$ws.Range("B4").Value = 30
$ws.Range("C4").Value = 30

# Row 5: was empty -> now Other Buttons/35/35
$ws.Range("A5").Value = "Other Buttons"
$ws.Range("B5").Value = 35
$ws.Range("C5").Value = 35

# Row 6: was empty -> now Slotboard/298/98
$ws.Range("A6").Value = "Slotboard"
$ws.Range("B6").Value = 298
$ws.Range("C6").Value = 98

# Row 7: was empty -> now DisplayTextBox/0/0
$ws.Range("A7").Value = "DisplayTextBox"
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0

# Update the view's active selection and scroll position to match the
# author's saved window state (cell C5 selected, scrolled so row 10 is
# the top visible row).
$ws.Range("C5").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
